# Updates cryptos list values per latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.748.04'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.59%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.324.54'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.04%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '579.74'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.47%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '175.11'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -4.18%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.588'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.46%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.320.96'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.21%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.178'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.77%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.576'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.77%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '45.38'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.19%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -2.19%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '657.99'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +4.04%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.860.85'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.27%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '8.40'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.86%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '67.627.77'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.88%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.00%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.315.92'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.33'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -2.01%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.886'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.91%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.35'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +5.34%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '16.95'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -3.91%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '98.39'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.50%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -3.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.66'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -4.09%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.23'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -3.60%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +2.83%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.41'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.26%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.22'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +7.52%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '567.28'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -6.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '10.92'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.45%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.105'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.11%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.07%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.662.95'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -7.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '56.23'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.55%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -6.90%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '34.43'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +5.39%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.24%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -2.79%  '
$ws.Range('B42').Value = 'ApeXProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.38'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.11'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -5.14%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.78%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0₃0661'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -3.86%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0404'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.59%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.25%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.06%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.29%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.08%  '
$ws.Range('B51').Value = 'CoreDAO'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.68'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +10.15%  '
